$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.902352809906006
$ws.Range("B1").Value = 2.084931135177612
$ws.Range("C1").Value = 2.432700872421265
$ws.Range("D1").Value = 3.758384943008423
$ws.Range("E1").Value = 1.258178949356079
